$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (avoid Excel auto-converting
# numeric-looking strings like "561.30" or "0.0000279" into numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.602.27'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.397.72'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.30'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.08'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.392.65'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.03%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.46%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000279'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.23'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.938.32'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.32'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.389.95'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '65.706.15'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.89'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '464.11'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.98'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.21%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +6.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.39'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.93'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.74'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.76'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.09'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.04%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '63.44'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.51'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '578.73'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.72%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.99'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.19%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0742'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.121.82'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.81'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.93%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.41%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.63%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.18'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.998'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.54'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.36%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.34%  '
